$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data in column A (record id) and column H (trailing id) are
# text-like numeric strings in the original file (shared-string cells),
# not numbers. Pre-format the ranges as Text so the new values are stored
# the same way (t="s") instead of being auto-coerced to numeric cells.
$ws.Range("A2:A16").NumberFormat = "@"
$ws.Range("H2:H16").NumberFormat = "@"

$ws.Range("A2").Value = "416489317"
$ws.Range("H2").Value = "8779"

$ws.Range("A3").Value = "941930440"
$ws.Range("H3").Value = "6618"

$ws.Range("A4").Value = "201335124"
$ws.Range("H4").Value = "8038"

$ws.Range("A5").Value = "792132756"
$ws.Range("H5").Value = "7022"

$ws.Range("A6").Value = "428986907"
$ws.Range("H6").Value = "1568"

$ws.Range("A7").Value = "267788365"
$ws.Range("H7").Value = "1041"

$ws.Range("A8").Value = "333545900"
$ws.Range("H8").Value = "8475"

$ws.Range("A9").Value = "900123606"
$ws.Range("H9").Value = "9471"

$ws.Range("A10").Value = "836455326"
$ws.Range("H10").Value = "9116"

$ws.Range("A11").Value = "805445054"
$ws.Range("H11").Value = "1286"

$ws.Range("A12").Value = "212848898"
$ws.Range("H12").Value = "7843"

$ws.Range("A13").Value = "767236544"
$ws.Range("H13").Value = "1863"

$ws.Range("A14").Value = "783495863"
$ws.Range("H14").Value = "7114"

$ws.Range("A15").Value = "394541246"
$ws.Range("H15").Value = "6758"

$ws.Range("A16").Value = "435238503"
$ws.Range("H16").Value = "4738"
